$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 1.732128977775574
$ws.Range("B1").Value = 2.589139699935913
$ws.Range("C1").Value = 2.740685939788818
$ws.Range("D1").Value = 3.09007453918457
$ws.Range("E1").Value = 3.353185892105103
